# Update cryptos list - GitHub Actions style data refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper to force a cell's value to remain plain text even when it looks
# numeric (e.g. "575.50", "63.00"), without permanently altering the
# cell's style/number format.
function Set-TextValue($range, $value) {
    $origStyle = $range.Style
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = $origStyle
}

# Row 2 - Bitcoin
Set-TextValue $ws.Range("D2") "60.214.03"
$ws.Range("E2").Value = "  +0.83%  "

# Row 3 - Ethereum
Set-TextValue $ws.Range("D3") "2.599.73"
$ws.Range("E3").Value = "  +0.05%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  -0.03%  "

# Row 5 - BNB
Set-TextValue $ws.Range("D5") "575.50"
$ws.Range("E5").Value = "  +3.38%  "

# Row 6 - Solana
Set-TextValue $ws.Range("D6") "142.86"
$ws.Range("E6").Value = "  +0.90%  "

# Row 7 - USDC
$ws.Range("E7").Value = "  +0.15%  "

# Row 8 - XRP
$ws.Range("E8").Value = "  -0.64%  "

# Row 9 - LidoStakedEther
Set-TextValue $ws.Range("D9") "2.604.18"
$ws.Range("E9").Value = "  -0.57%  "

# Row 10 - Toncoin
Set-TextValue $ws.Range("D10") "6.57"
$ws.Range("E10").Value = "  -1.87%  "

# Row 11 - Dogecoin
$ws.Range("E11").Value = "  +0.81%  "

# Row 12 - TRON
$ws.Range("E12").Value = "  -3.50%  "

# Row 13 - Cardano
$ws.Range("E13").Value = "  +0.33%  "

# Row 14 - WrappedliquidstakedEther2.0
Set-TextValue $ws.Range("D14") "3.057.24"

# Row 15 - Avalanche
Set-TextValue $ws.Range("D15") "24.35"
$ws.Range("E15").Value = "  +3.80%  "

# Row 16 - WrappedBTC
Set-TextValue $ws.Range("D16") "60.217.69"
$ws.Range("E16").Value = "  +0.87%  "

# Row 17 - ShibaInu
$ws.Range("E17").Value = "  +2.01%  "

# Row 18 - WrappedEther
Set-TextValue $ws.Range("D18") "2.605.86"
$ws.Range("E18").Value = "  -0.28%  "

# Row 19 - Chainlink
Set-TextValue $ws.Range("D19") "11.36"
$ws.Range("E19").Value = "  +7.04%  "

# Row 20 - Polkadot
$ws.Range("E20").Value = "  +0.15%  "

# Row 21 - BitcoinCash
Set-TextValue $ws.Range("D21") "346.19"
$ws.Range("E21").Value = "  +1.15%  "

# Row 22 - Uniswap
Set-TextValue $ws.Range("D22") "6.89"
$ws.Range("E22").Value = "  +1.87%  "

# Row 23 - Dai
$ws.Range("E23").Value = "  -0.02%  "

# Row 24 - Polygon
Set-TextValue $ws.Range("D24") "0.530"
$ws.Range("E24").Value = "  +2.17%  "

# Row 25 - Litecoin
Set-TextValue $ws.Range("D25") "63.00"
$ws.Range("E25").Value = "  +0.79%  "

# Row 26 - Binance-PegBSC-USD
$ws.Range("E26").Value = "  +0.15%  "

# Row 27 - Kaspa
Set-TextValue $ws.Range("D27") "0.159"
$ws.Range("E27").Value = "  -0.21%  "

# Row 28 - InternetComputer(DFINITY)
Set-TextValue $ws.Range("D28") "8.01"
$ws.Range("E28").Value = "  +5.68%  "

# Row 29 - PEPE
$ws.Range("E29").Value = "  +2.02%  "

# Row 30 - PancakeSwap
$ws.Range("E30").Value = "  +9.64%  "

# Row 31 - Aptos
Set-TextValue $ws.Range("D31") "6.38"
$ws.Range("E31").Value = "  +3.54%  "

# Row 33 - Monero
Set-TextValue $ws.Range("D33") "166.46"
$ws.Range("E33").Value = "  +4.94%  "

# Row 34 - EthereumClassic
Set-TextValue $ws.Range("D34") "19.41"
$ws.Range("E34").Value = "  +0.09%  "

# Row 35 - was NEARProtocol, now ImmutableX
$ws.Range("B35").Value = "ImmutableX"
$ws.Range("C35").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
Set-TextValue $ws.Range("D35") "1.30"
$ws.Range("E35").Value = "  +9.66%  "

# Row 36 - was ImmutableX, now NEARProtocol
$ws.Range("B36").Value = "NEARProtocol"
$ws.Range("C36").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
Set-TextValue $ws.Range("D36") "4.27"
$ws.Range("E36").Value = "  +3.38%  "

# Row 37 - Fetch.AI
Set-TextValue $ws.Range("D37") "0.985"
$ws.Range("E37").Value = "  +7.47%  "

# Row 38 - Stacks
$ws.Range("E38").Value = "  +6.70%  "

# Row 39 - OKB
Set-TextValue $ws.Range("D39") "38.06"
$ws.Range("E39").Value = "  +0.82%  "

# Row 40 - Bittensor
Set-TextValue $ws.Range("D40") "312.93"
$ws.Range("E40").Value = "  +7.03%  "

# Row 41 - Filecoin
$ws.Range("E41").Value = "  +4.58%  "

# Row 42 - SuiNetwork
Set-TextValue $ws.Range("D42") "0.837"
$ws.Range("E42").Value = "  -0.80%  "

# Row 43 - Aave
Set-TextValue $ws.Range("D43") "135.05"
$ws.Range("E43").Value = "  -2.87%  "

# Row 44 - Stellar
Set-TextValue $ws.Range("D44") "0.0993"
$ws.Range("E44").Value = "  +1.35%  "

# Row 45 - FirstDigitalUSD
$ws.Range("E45").Value = "  +0.15%  "

# Row 46 - EnergySwap
Set-TextValue $ws.Range("D46") "19.84"
$ws.Range("E46").Value = "  +2.70%  "

# Row 47 - Hedera
Set-TextValue $ws.Range("D47") "0.0551"
$ws.Range("E47").Value = "  +2.07%  "

# Row 48 - was RenderToken, now Mantle
$ws.Range("B48").Value = "Mantle"
$ws.Range("C48").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
Set-TextValue $ws.Range("D48") "0.604"
$ws.Range("E48").Value = "  +0.50%  "

# Row 49 - was Mantle, now RenderToken
$ws.Range("B49").Value = "RenderToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextValue $ws.Range("D49") "4.99"
$ws.Range("E49").Value = "  +4.05%  "

# Row 50 - was InjectiveProtocol, now VeChain
$ws.Range("B50").Value = "VeChain"
$ws.Range("C50").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextValue $ws.Range("D50") "0.0241"
$ws.Range("E50").Value = "  +0.37%  "

# Row 51 - was VeChain, now InjectiveProtocol
$ws.Range("B51").Value = "InjectiveProtocol"
$ws.Range("C51").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
Set-TextValue $ws.Range("D51") "19.94"
$ws.Range("E51").Value = "  +5.18%  "
